$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5-16 hold pairs of ITC experiments for guest01..guest06: a
# "buffer into guestNN" row (SamplePrepMethod = Chodera Load Cell
# Without Cleaning Cell After.setup) immediately followed by a
# "host into guestNN" row (SamplePrepMethod = Plates Quick.setup).
#
# This reverses the order of each pair: the Host2Guest info
# (SampleName/SamplePrepMethod/PipetteConcentration) moves up to the
# first row of the pair and the Buff2Guest info moves down to the
# second row. The CellConcentration (column F) value that belonged to
# the (second, Host2Guest) row is kept by both rows of the pair, since
# it is tied to the underlying data file rather than to the swapped
# sample-prep ordering.

$rows = @(5, 7, 9, 11, 13, 15)

foreach ($r1 in $rows) {
    $r2 = $r1 + 1

    # Capture current (pre-swap) values for the pair of rows.
    $b1 = $ws.Cells.Item($r1, 2).Value()   # SampleName (Buff2Guest)
    $c1 = $ws.Cells.Item($r1, 3).Value()   # SamplePrepMethod (Buff2Guest)
    $g1 = $ws.Cells.Item($r1, 7).Value()   # PipetteConcentration (Buff2Guest)

    $b2 = $ws.Cells.Item($r2, 2).Value()   # SampleName (Host2Guest)
    $c2 = $ws.Cells.Item($r2, 3).Value()   # SamplePrepMethod (Host2Guest)
    $f2 = $ws.Cells.Item($r2, 6).Value()   # CellConcentration (Host2Guest)
    $g2 = $ws.Cells.Item($r2, 7).Value()   # PipetteConcentration (Host2Guest)

    # First row of the pair becomes the Host2Guest entry.
    $ws.Cells.Item($r1, 2).Value = $b2
    $ws.Cells.Item($r1, 3).Value = $c2
    $ws.Cells.Item($r1, 6).Value = $f2
    $ws.Cells.Item($r1, 7).Value = $g2

    # Second row of the pair becomes the Buff2Guest entry, but keeps
    # the CellConcentration value from the Host2Guest data file.
    $ws.Cells.Item($r2, 2).Value = $b1
    $ws.Cells.Item($r2, 3).Value = $c1
    $ws.Cells.Item($r2, 6).Value = $f2
    $ws.Cells.Item($r2, 7).Value = $g1
}
